# add preview in datasets
# Adds two new rows (dataset.preview.table / dataset.preview.line) with
# DataWatch-style query formulas right after the "dataset.commit.id" header
# row on the "metadata" sheet, matching the existing look (wrapped text,
# vertically centered, 120pt tall rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("metadata")

# Make room for the two new rows (old row 4 and everything below shifts
# down by two).
$ws.Range("A4:A5").EntireRow.Insert()

$previewTableFormula = @"
source(ds:'{{dataset.id}}');
query([
  { dim:'time', role:'row', items:[] },
  { dim:'indicator', role:'col', items:[] } 
]);
format(p:3);
order(dir:'row', index:-1, asc:'az');
limit(start:0, length:5);
"@

$previewLineFormula = @"
source(ds:'{{dataset.id}}');
query([
  { dim:'time', role:'row', items:[] },
  { dim:'indicator', role:'col', items:[] } 
]);
format(p:3);
order(dir:'row', index:-1, asc:'az');
line(x:-1);
"@

$ws.Cells.Item(4, 1).Value() = "dataset.preview.table"
$ws.Cells.Item(4, 2).Value() = $previewTableFormula
$ws.Cells.Item(5, 1).Value() = "dataset.preview.line"
$ws.Cells.Item(5, 2).Value() = $previewLineFormula

# Match the formatting of the rest of the key/value table, but wrap text
# and grow the rows so the multi-line formulas are fully visible.
$ws.Range("A4:B5").Style = $ws.Range("A3:B3").Style
$ws.Range("A4:B5").VerticalAlignment = -4108
$ws.Range("A4:B5").WrapText = 1
$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 120

$ws.Range("B6").Select()
